$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation (2026/01/10 土, hour 15, ranking 201) was recorded and
# inserted as row 620, pushing every subsequent row down by one (old row
# 661 -> new row 662). Use a real row insert so the shift happens exactly
# like it would from the Excel UI / COM automation.
$ws.Rows.Item(620).Insert()

# Column A holds the date as literal text (e.g. "2026/01/10"), matching
# every other row in the sheet (t="inlineStr"/shared string, not a real
# date serial). Mark the cell as Text *before* assigning the value so
# Excel's autodetect doesn't silently convert the "YYYY/MM/DD"-looking
# string into a date serial number.
$ws.Cells.Item(620, 1).NumberFormat = "@"
$ws.Cells.Item(620, 1).Value = "2026/01/10"
$ws.Cells.Item(620, 2).Value = "土"
$ws.Cells.Item(620, 3).Value = 15
$ws.Cells.Item(620, 4).Value = 201

# Forcing the "@" text format above stamps A620 with its own new style
# (quote-prefix / text numFmt) that none of the surrounding data cells
# carry. Re-sync its formatting with the untouched cell directly below it
# (same column, same intended "plain text date" look) so the inserted row
# ends up styled identically to every other data row.
$ws.Cells.Item(621, 1).Copy()
$ws.Cells.Item(620, 1).PasteSpecial(-4122)
